# Collapse the wrapped, multi-line `aggregate(...)` and `rollapplyr(...)`
# calls in the "Given Code" answer-key block back onto single lines: each
# continuation that used to be "<text>,<manual line break><indent>" becomes
# "<text>, " (a manual line break + the following indentation run replaced
# by a single trailing space). Only the blank "___" / "filled_co2" answer
# key occurrences (23-space indent under `aggregate`, 19-space indent under
# `rollapplyr`) are touched; the earlier "Outline" walkthrough code with the
# filled-in answers and 2-space indents is left untouched.

$d = $word.ActiveDocument

# Manual line break character, as produced by <w:br/> in Find/Replace text.
$lineBreak = [char]11

function Build-Spaces($count) {
    $s = ""
    for ($i = 0; $i -lt $count; $i++) {
        $s += " "
    }
    return $s
}

$indent23 = Build-Spaces 23
$indent19 = Build-Spaces 19

# 1) "aggregate(___," + <br/> + 23 spaces  ->  "aggregate(___, "
$old1 = "(___," + $lineBreak + $indent23
$new1 = "(___, "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2) " ___," (the `by = ___,` continuation) + <br/> + 23 spaces -> " ___, "
$old2 = " ___," + $lineBreak + $indent23
$new2 = " ___, "
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) "rollapplyr(filled_co2," + <br/> + 19 spaces -> "rollapplyr(filled_co2, "
$old3 = "(filled_co2," + $lineBreak + $indent19
$new3 = "(filled_co2, "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 4) " ___," (the `FUN = ___,` continuation) + <br/> + 19 spaces -> " ___, "
$old4 = " ___," + $lineBreak + $indent19
$new4 = " ___, "
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

Write-Host "Done collapsing aggregate()/rollapplyr() continuation lines."
